# Updated cryptos list (price/volume refresh) - rows 2..51 on the active sheet.
# Numeric-looking price strings (single decimal point) are written with a
# leading apostrophe so Excel keeps them as literal text instead of coercing
# them into floating-point numbers (matching the source data's inline-string
# formatting, e.g. "314.73" rather than 314.73000000000002).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.587.90"
$ws.Cells.Item(2, 5).Value = "  +2.42%  "

$ws.Cells.Item(3, 4).Value = "1.791.48"
$ws.Cells.Item(3, 5).Value = "  +4.25%  "

$ws.Cells.Item(4, 4).Value = "'0.9997"
$ws.Cells.Item(4, 5).Value = "  -0.13%  "

$ws.Cells.Item(5, 4).Value = "'314.73"

$ws.Cells.Item(6, 4).Value = "'0.9995"

$ws.Cells.Item(7, 4).Value = "'0.5354"
$ws.Cells.Item(7, 5).Value = "  +10.01%  "

$ws.Cells.Item(8, 4).Value = "'0.3782"
$ws.Cells.Item(8, 5).Value = "  +8.14%  "

$ws.Cells.Item(9, 4).Value = "'43.05"
$ws.Cells.Item(9, 5).Value = "  +2.66%  "

$ws.Cells.Item(10, 4).Value = "'0.07539"
$ws.Cells.Item(10, 5).Value = "  +4.01%  "

$ws.Cells.Item(11, 4).Value = "'1.115"
$ws.Cells.Item(11, 5).Value = "  +6.68%  "

$ws.Cells.Item(12, 4).Value = "'0.9997"
$ws.Cells.Item(12, 5).Value = "  -0.11%  "

$ws.Cells.Item(13, 4).Value = "'20.94"
$ws.Cells.Item(13, 5).Value = "  +5.35%  "

$ws.Cells.Item(14, 4).Value = "'6.193"
$ws.Cells.Item(14, 5).Value = "  +5.91%  "

$ws.Cells.Item(15, 4).Value = "1.789.33"
$ws.Cells.Item(15, 5).Value = "  +4.35%  "

$ws.Cells.Item(16, 4).Value = "'7.092"
$ws.Cells.Item(16, 5).Value = "  +3.50%  "

$ws.Cells.Item(17, 4).Value = "'90.99"
$ws.Cells.Item(17, 5).Value = "  +5.13%  "

$ws.Cells.Item(18, 4).Value = "'0.00001068"
$ws.Cells.Item(18, 5).Value = "  +2.92%  "

$ws.Cells.Item(19, 4).Value = "'0.06503"
$ws.Cells.Item(19, 5).Value = "  +2.19%  "

$ws.Cells.Item(20, 4).Value = "'0.9992"
$ws.Cells.Item(20, 5).Value = "  -0.06%  "

$ws.Cells.Item(21, 4).Value = "'16.98"
$ws.Cells.Item(21, 5).Value = "  +2.98%  "

$ws.Cells.Item(22, 4).Value = "'5.946"
$ws.Cells.Item(22, 5).Value = "  +5.34%  "

$ws.Cells.Item(23, 4).Value = "27.608.23"
$ws.Cells.Item(23, 5).Value = "  +2.27%  "

$ws.Cells.Item(25, 4).Value = "'2.096"
$ws.Cells.Item(25, 5).Value = "  +0.71%  "

$ws.Cells.Item(26, 4).Value = "'20.54"
$ws.Cells.Item(26, 5).Value = "  +2.72%  "

$ws.Cells.Item(27, 4).Value = "'155.56"
$ws.Cells.Item(27, 5).Value = "  +1.32%  "

$ws.Cells.Item(28, 4).Value = "'2.394"
$ws.Cells.Item(28, 5).Value = "  +15.32%  "

$ws.Cells.Item(29, 4).Value = "1.994.10"
$ws.Cells.Item(29, 5).Value = "  +4.24%  "

$ws.Cells.Item(30, 4).Value = "'122.21"
$ws.Cells.Item(30, 5).Value = "  +1.05%  "

$ws.Cells.Item(31, 4).Value = "'1.118"
$ws.Cells.Item(31, 5).Value = "  +8.80%  "

$ws.Cells.Item(32, 4).Value = "'0.1031"
$ws.Cells.Item(32, 5).Value = "  +10.94%  "

$ws.Cells.Item(33, 4).Value = "'5.705"
$ws.Cells.Item(33, 5).Value = "  +6.70%  "

$ws.Cells.Item(34, 4).Value = "'3.624"
$ws.Cells.Item(34, 5).Value = "  +1.20%  "

$ws.Cells.Item(35, 4).Value = "'0.02290"
$ws.Cells.Item(35, 5).Value = "  +5.04%  "

$ws.Cells.Item(36, 4).Value = "'8.658"
$ws.Cells.Item(36, 5).Value = "  +15.29%  "

$ws.Cells.Item(37, 4).Value = "'0.06030"
$ws.Cells.Item(37, 5).Value = "  +2.57%  "

$ws.Cells.Item(38, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(38, 4).Value = "'4.996"
$ws.Cells.Item(38, 5).Value = "  +5.61%  "

$ws.Cells.Item(39, 2).Value = "Aptos"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(39, 4).Value = "'11.44"
$ws.Cells.Item(39, 5).Value = "  +4.42%  "

$ws.Cells.Item(40, 4).Value = "'0.2086"
$ws.Cells.Item(40, 5).Value = "  +4.67%  "

$ws.Cells.Item(41, 4).Value = "'0.6256"
$ws.Cells.Item(41, 5).Value = "  +4.03%  "

$ws.Cells.Item(42, 4).Value = "'1.411"
$ws.Cells.Item(42, 5).Value = "  -2.87%  "

$ws.Cells.Item(43, 4).Value = "'0.9991"
$ws.Cells.Item(43, 5).Value = "  +0.01%  "

$ws.Cells.Item(44, 4).Value = "'1.146"
$ws.Cells.Item(44, 5).Value = "  +4.65%  "

$ws.Cells.Item(45, 4).Value = "'13.36"
$ws.Cells.Item(45, 5).Value = "  +5.19%  "

$ws.Cells.Item(46, 4).Value = "'0.5870"
$ws.Cells.Item(46, 5).Value = "  +4.24%  "

$ws.Cells.Item(47, 4).Value = "'3.637"
$ws.Cells.Item(47, 5).Value = "  +1.80%  "

$ws.Cells.Item(48, 4).Value = "'121.91"
$ws.Cells.Item(48, 5).Value = "  +2.96%  "

$ws.Cells.Item(49, 4).Value = "'1.921"
$ws.Cells.Item(49, 5).Value = "  +4.85%  "

$ws.Cells.Item(50, 4).Value = "'1.133"
$ws.Cells.Item(50, 5).Value = "  +1.97%  "

$ws.Cells.Item(51, 5).Value = "  +1.48%  "
